$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing runtime values (days 1-8) with re-measured timings
$ws.Range("B3").Value = 0.0027379599999999998
$ws.Range("B4").Value = 0.00357888
$ws.Range("B5").Value = 0.00200994
$ws.Range("B6").Value = 0.0043233999999999998
$ws.Range("B7").Value = 0.010369079999999999
$ws.Range("B8").Value = 15.880498640000001
$ws.Range("B9").Value = 2.1064885200000001
$ws.Range("B10").Value = 0.00309516

# Add new rows for days 9-12
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 0.14254684000000001

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 0.015152519999999999

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 0.090809299999999996

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 0.064846819999999999

# Update the selection to match the new data extent
$ws.Range("A3:B14").Select()
